$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Git section
$ws.Range("B10").Value = "Git"
$ws.Range("C11").Value = "git remote set-url origin https://<githubtoken>@github.com/<username>/<repositoryname>.git"
$ws.Range("C13").Value = "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"
$ws.Range("C12").Value = "git remote set-url origin https://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@github.com/PCnslt/Python.git"

# Python section
$ws.Range("B15").Value = "Python2 and 3"
$ws.Range("C16").Value = "Use PyCharm as IDE"

# MAC Address Changer section
$ws.Range("A18").Value = "MAC Address Changer"
$ws.Range("B19").Value = "What is MAC"
$ws.Range("C19").Value = "Media Access Control"

$ws.Range("B21").Value = "Change MAC address using terminal"
$ws.Range("C22").Value = "ifconfig"
$ws.Range("C23").Value = "ifconfig wlan0 down"
$ws.Range("C24").Value = "ifconfig wlan0 hw ether 00:11:22:33:44:55"
$ws.Range("C25").Value = "ifconfig wlan0 up"
$ws.Range("C26").Value = "ifconfig wlan0"

$ws.Range("B28").Value = "User input"
$ws.Range("C29").Value = "https://docs.python.org/2/library/functions.html#raw_input"

$ws.Range("B31").Value = "RegEx"
$ws.Range("C32").Value = "pythex.org"

# Update view state to match target: window geometry, scrolled to A17, selection C33
$win = $excel.ActiveWindow
$win.Top = 2070
$win.Left = 3720
$win.Width = 21600
$win.Height = 11295
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("C33").Select()
